$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For rows 1..64: convert the numeric "division" value in column D (1, 2 or 3)
# into the text labels "S1"/"S2"/"S3", and add a new column E with the
# constant label "A" (e.g. department/year marker), for each data row.
for ($r = 1; $r -le 64; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 4).Value = "S" + [string]$d
    $ws.Cells.Item($r, 5).Value = "A"
}

# Update the active selection to match the saved view state (single cell F5,
# scrolled back to the top-left of the sheet).
$ws.Range("F5").Select()
